{"js": "const replacements = [\n  [\"2024-03-02 Saturday\", \"2024-03-03 Sunday\"],\n  [\"18\u00d791=\", \"87\u00d721=\"],\n  [\"78\u00d758=\", \"56\u00d772=\"],\n  [\"18\u00d790=\", \"14\u00d722=\"],\n  [\"42\u00d797=\", \"40\u00d775=\"],\n  [\"49\u00d777=\", \"66\u00d715=\"],\n  [\"62\u00d750=\", \"56\u00d735=\"],\n  [\"21\u00d797=\", \"24\u00d767=\"],\n  [\"18\u00d779=\", \"59\u00d713=\"],\n  [\"84\u00d723=\", \"91\u00d734=\"],\n  [\"86\u00d712=\", \"72\u00d751=\"],\n  [\"85\u00d734=\", \"64\u00d754=\"],\n  [\"24\u00d721=\", \"34\u00d765=\"],\n  [\"50\u00d756=\", \"82\u00d732=\"],\n  [\"38\u00d788=\", \"13\u00d779=\"],\n  [\"64\u00d756=\", \"31\u00d786=\"],\n  [\"76\u00d711=\", \"37\u00d714=\"],\n  [\"47\u00d754=\", \"23\u00d793=\"],\n  [\"68\u00d755=\", \"44\u00d775=\"],\n  [\"58\u00d781=\", \"33\u00d747=\"],\n  [\"37\u00d722=\", \"52\u00d750=\"],\n  [\"86\u00d729=\", \"20\u00d778=\"],\n  [\"13\u00d763=\", \"23\u00d787=\"],\n  [\"69\u00d725=\", \"77\u00d755=\"],\n  [\"77\u00d725=\", \"97\u00d722=\"],\n  [\"66\u00d755=\", \"54\u00d767=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-03-02 Saturday\", \"2024-03-03 Sunday\"),\n    @(\"18\u00d791=\", \"87\u00d721=\"),\n    @(\"78\u00d758=\", \"56\u00d772=\"),\n    @(\"18\u00d790=\", \"14\u00d722=\"),\n    @(\"42\u00d797=\", \"40\u00d775=\"),\n    @(\"49\u00d777=\", \"66\u00d715=\"),\n    @(\"62\u00d750=\", \"56\u00d735=\"),\n    @(\"21\u00d797=\", \"24\u00d767=\"),\n    @(\"18\u00d779=\", \"59\u00d713=\"),\n    @(\"84\u00d723=\", \"91\u00d734=\"),\n    @(\"86\u00d712=\", \"72\u00d751=\"),\n    @(\"85\u00d734=\", \"64\u00d754=\"),\n    @(\"24\u00d721=\", \"34\u00d765=\"),\n    @(\"50\u00d756=\", \"82\u00d732=\"),\n    @(\"38\u00d788=\", \"13\u00d779=\"),\n    @(\"64\u00d756=\", \"31\u00d786=\"),\n    @(\"76\u00d711=\", \"37\u00d714=\"),\n    @(\"47\u00d754=\", \"23\u00d793=\"),\n    @(\"68\u00d755=\", \"44\u00d775=\"),\n    @(\"58\u00d781=\", \"33\u00d747=\"),\n    @(\"37\u00d722=\", \"52\u00d750=\"),\n    @(\"86\u00d729=\", \"20\u00d778=\"),\n    @(\"13\u00d763=\", \"23\u00d787=\"),\n    @(\"69\u00d725=\", \"77\u00d755=\"),\n    @(\"77\u00d725=\", \"97\u00d722=\"),\n    @(\"66\u00d755=\", \"54\u00d767=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\n"}
